$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1010.64703
$ws.Range("I28").Value = 760.9375
$ws.Range("J28").Value = 5006
$ws.Range("K28").Value = 760.9375
$ws.Range("L28").Value = 5006
$ws.Range("M28").Value = -275.9375
$ws.Range("N28").Value = -5976
$ws.Range("H38").Value = 436.27274
$ws.Range("I38").Value = 436.27274
$ws.Range("K38").Value = 1308.81822
$ws.Range("M38").Value = -936.8182200000001
$ws.Range("H61").Value = 356.25
$ws.Range("I61").Value = 356.25
$ws.Range("K61").Value = 1068.75
$ws.Range("M61").Value = -896.75
$ws.Range("H62").Value = 4754.85
$ws.Range("I62").Value = 4384.5
$ws.Range("K62").Value = 4384.5
$ws.Range("M62").Value = -3760.5
$ws.Range("H63").Value = 25271
$ws.Range("J63").Value = 25271
$ws.Range("L63").Value = 25271
$ws.Range("N63").Value = -26519
$ws.Range("H65").Value = 4754.85
$ws.Range("I65").Value = 4384.5
$ws.Range("K65").Value = 21922.5
$ws.Range("M65").Value = -18802.5
$ws.Range("H66").Value = 25271
$ws.Range("J66").Value = 25271
$ws.Range("L66").Value = 75813
$ws.Range("N66").Value = -82053
$ws.Range("H92").Value = 235
$ws.Range("I92").Value = 243.21739
$ws.Range("K92").Value = 243.21739
$ws.Range("M92").Value = 1004.78261
$ws.Range("H99").Value = 111130260
$ws.Range("J99").Value = 333334180
$ws.Range("L99").Value = 1000002540
$ws.Range("N99").Value = -1000005536
$ws.Range("H100").Value = 6473.2856
$ws.Range("J100").Value = 4504.4546
$ws.Range("L100").Value = 4504.4546
$ws.Range("N100").Value = -5586.4546
$ws.Range("H101").Value = 36500932
$ws.Range("I101").Value = 910195.8
$ws.Range("K101").Value = 2730587.4
$ws.Range("M101").Value = -2728965.4
$ws.Range("H118").Value = 251025
$ws.Range("I118").Value = 251025
$ws.Range("K118").Value = 753075
$ws.Range("M118").Value = -751418
$ws.Range("H138").Value = 2528.742
$ws.Range("I138").Value = 2684.077
$ws.Range("J138").Value = 2416.5557
$ws.Range("K138").Value = 8052.231000000001
$ws.Range("L138").Value = 7249.6671
$ws.Range("M138").Value = -2912.231000000001
$ws.Range("N138").Value = -17529.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1890.5652
$ws.Range("I2").Value = 1887.5555
$ws.Range("K2").Value = 1887.5555
$ws.Range("M2").Value = -1774.5555
$ws.Range("H28").Value = 15499.125
$ws.Range("I28").Value = 13353.429
$ws.Range("K28").Value = 13353.429
$ws.Range("M28").Value = -13161.429
$ws.Range("H32").Value = 2400.5217
$ws.Range("I32").Value = 2400.5217
$ws.Range("K32").Value = 2400.5217
$ws.Range("M32").Value = -2113.5217
$ws.Range("H61").Value = 2137.6667
$ws.Range("I61").Value = 2137.6667
$ws.Range("K61").Value = 2137.6667
$ws.Range("M61").Value = -1925.6667
$ws.Range("H63").Value = 3739.25
$ws.Range("I63").Value = 3739.25
$ws.Range("K63").Value = 3739.25
$ws.Range("M63").Value = -3053.25
$ws.Range("H66").Value = 3739.25
$ws.Range("I66").Value = 3739.25
$ws.Range("K66").Value = 18696.25
$ws.Range("M66").Value = -15264.25
$ws.Range("H99").Value = 15499.125
$ws.Range("I99").Value = 13353.429
$ws.Range("K99").Value = 13353.429
$ws.Range("M99").Value = -10358.429
$ws.Range("H116").Value = 1890.5652
$ws.Range("I116").Value = 1887.5555
$ws.Range("K116").Value = 1887.5555
$ws.Range("M116").Value = 406.4445000000001
$ws.Range("H132").Value = 4284.143
$ws.Range("I132").Value = 3831.5
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 11494.5
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -8964.5
$ws.Range("N132").Value = -26060
$ws.Range("H136").Value = 2137.6667
$ws.Range("I136").Value = 2137.6667
$ws.Range("K136").Value = 6413.000100000001
$ws.Range("M136").Value = -3863.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1890.5652
$ws.Range("I3").Value = 1887.5555
$ws.Range("K3").Value = 1887.5555
$ws.Range("M3").Value = -1773.5555
$ws.Range("H20").Value = 9197.84
$ws.Range("I20").Value = 9325.950000000001
$ws.Range("K20").Value = 9325.950000000001
$ws.Range("M20").Value = -9078.950000000001
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H82").Value = 22099.857
$ws.Range("H85").Value = 22099.857
$ws.Range("H86").Value = 3284.4285
$ws.Range("J86").Value = 707.25
$ws.Range("L86").Value = 707.25
$ws.Range("N86").Value = -2953.25
$ws.Range("H89").Value = 3284.4285
$ws.Range("J89").Value = 707.25
$ws.Range("L89").Value = 3536.25
$ws.Range("N89").Value = -14768.25
$ws.Range("H94").Value = 1024.75
$ws.Range("I94").Value = 1024.75
$ws.Range("K94").Value = 1024.75
$ws.Range("M94").Value = -573.75
$ws.Range("H107").Value = 114156.22
$ws.Range("I107").Value = 169833.17
$ws.Range("J107").Value = 2802.3333
$ws.Range("K107").Value = 169833.17
$ws.Range("L107").Value = 2802.3333
$ws.Range("M107").Value = -167913.17
$ws.Range("N107").Value = -6642.3333
$ws.Range("H141").Value = 64639.75
$ws.Range("J141").Value = 64639.75
$ws.Range("L141").Value = 64639.75
$ws.Range("N141").Value = -74999.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17816.666
$ws.Range("I41").Value = 11725
$ws.Range("K41").Value = 11725
$ws.Range("M41").Value = -11297
$ws.Range("H107").Value = 3197.625
$ws.Range("I107").Value = 3105.182
$ws.Range("J107").Value = 3401
$ws.Range("K107").Value = 3105.182
$ws.Range("L107").Value = 3401
$ws.Range("M107").Value = -1185.182
$ws.Range("N107").Value = -7241
$ws.Range("H141").Value = 440261.84
$ws.Range("J141").Value = 440261.84
$ws.Range("L141").Value = 440261.84
$ws.Range("N141").Value = -450621.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1054.2
$ws.Range("I5").Value = 923.6667
$ws.Range("K5").Value = 2771.0001
$ws.Range("M5").Value = -2659.0001
$ws.Range("H8").Value = 828.375
$ws.Range("I8").Value = 828.375
$ws.Range("K8").Value = 2485.125
$ws.Range("M8").Value = -2346.125
$ws.Range("H131").Value = 1966.375
$ws.Range("J131").Value = 1919.4166
$ws.Range("L131").Value = 5758.2498
$ws.Range("N131").Value = -15838.2498
$ws.Range("H134").Value = 2367.4614
$ws.Range("I134").Value = 1731.4166
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 5194.2498
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -124.2497999999996
$ws.Range("N134").Value = -40140
$ws.Range("H135").Value = 1054.2
$ws.Range("I135").Value = 923.6667
$ws.Range("K135").Value = 8313.0003
$ws.Range("M135").Value = -5778.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3852.3076
$ws.Range("I113").Value = 4086.875
$ws.Range("K113").Value = 4086.875
$ws.Range("M113").Value = -1916.875
$ws.Range("H122").Value = 1201.8
$ws.Range("I122").Value = 1143.5
$ws.Range("K122").Value = 3430.5
$ws.Range("M122").Value = -980.5
$ws.Range("H123").Value = 40980.23
$ws.Range("I123").Value = 10295.25
$ws.Range("J123").Value = 54618
$ws.Range("K123").Value = 10295.25
$ws.Range("L123").Value = 54618
$ws.Range("M123").Value = -7845.25
$ws.Range("N123").Value = -59518
$ws.Range("H132").Value = 7449.9644
$ws.Range("I132").Value = 6307.15
$ws.Range("J132").Value = 10307
$ws.Range("K132").Value = 18921.45
$ws.Range("L132").Value = 30921
$ws.Range("M132").Value = -16391.45
$ws.Range("N132").Value = -35981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1935.3334
$ws.Range("I82").Value = 2015.6
$ws.Range("J82").Value = 1835
$ws.Range("K82").Value = 2015.6
$ws.Range("L82").Value = 1835
$ws.Range("M82").Value = -1654.6
$ws.Range("N82").Value = -2557
$ws.Range("H85").Value = 1935.3334
$ws.Range("I85").Value = 2015.6
$ws.Range("J85").Value = 1835
$ws.Range("K85").Value = 2015.6
$ws.Range("L85").Value = 1835
$ws.Range("M85").Value = -767.5999999999999
$ws.Range("N85").Value = -4331
$ws.Range("H132").Value = 3033.4412
$ws.Range("I132").Value = 2669.0967
$ws.Range("K132").Value = 8007.2901
$ws.Range("M132").Value = -5477.2901

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9999
$ws.Range("J69").Value = 9999
$ws.Range("L69").Value = 9999
$ws.Range("N69").Value = -11497
$ws.Range("H72").Value = 9999
$ws.Range("J72").Value = 9999
$ws.Range("L72").Value = 29997
$ws.Range("N72").Value = -37485
$ws.Range("H107").Value = 16667939
$ws.Range("I107").Value = 1383.95
$ws.Range("K107").Value = 4151.85
$ws.Range("M107").Value = -2231.85
$ws.Range("H132").Value = 6987.3477
$ws.Range("I132").Value = 7879.421
$ws.Range("K132").Value = 23638.263
$ws.Range("M132").Value = -21108.263
